# Auto-generated edit script: round prediction/continent statistics
$wb = $excel.ActiveWorkbook
$wsPred = $wb.Worksheets.Item("Prediction_statistics")
$wsCont = $wb.Worksheets.Item("Continent_statistics")

# --- Prediction_statistics: round Mean_1/2_ci (D) and Production (G) to 2 decimals ---
$wsPred.Cells.Item(2, 4).Value = 12.71
$wsPred.Cells.Item(2, 7).Value = 909268060.27
$wsPred.Cells.Item(3, 4).Value = 9.41
$wsPred.Cells.Item(3, 7).Value = 910587344.24
$wsPred.Cells.Item(4, 4).Value = 8.03
$wsPred.Cells.Item(4, 7).Value = 671495828.91
$wsPred.Cells.Item(5, 4).Value = 4.94
$wsPred.Cells.Item(5, 7).Value = 509178836.77
$wsPred.Cells.Item(6, 4).Value = 0
$wsPred.Cells.Item(7, 4).Value = 0
$wsPred.Cells.Item(8, 4).Value = 12.98
$wsPred.Cells.Item(8, 7).Value = 381201572.84
$wsPred.Cells.Item(9, 4).Value = 10.24
$wsPred.Cells.Item(9, 7).Value = 384223628.3
$wsPred.Cells.Item(10, 4).Value = 6.72
$wsPred.Cells.Item(10, 7).Value = 263126743.8
$wsPred.Cells.Item(11, 4).Value = 5.85
$wsPred.Cells.Item(11, 7).Value = 247565762.45
$wsPred.Cells.Item(12, 4).Value = 0
$wsPred.Cells.Item(13, 4).Value = 0
$wsPred.Cells.Item(14, 4).Value = 5.58
$wsPred.Cells.Item(14, 7).Value = 177505211.59
$wsPred.Cells.Item(15, 4).Value = 3.21
$wsPred.Cells.Item(15, 7).Value = 177557093.87
$wsPred.Cells.Item(16, 4).Value = 4.11
$wsPred.Cells.Item(16, 7).Value = 152081931.41
$wsPred.Cells.Item(17, 4).Value = 2.02
$wsPred.Cells.Item(17, 7).Value = 111003266.03
$wsPred.Cells.Item(18, 4).Value = 0
$wsPred.Cells.Item(19, 4).Value = 0
$wsPred.Cells.Item(20, 4).Value = 6.8
$wsPred.Cells.Item(20, 7).Value = 617000745.92
$wsPred.Cells.Item(21, 4).Value = 4.5
$wsPred.Cells.Item(21, 7).Value = 617912194.77
$wsPred.Cells.Item(22, 4).Value = 2.71
$wsPred.Cells.Item(22, 7).Value = 407583172.81
$wsPred.Cells.Item(23, 4).Value = 1.89
$wsPred.Cells.Item(23, 7).Value = 328468517.69
$wsPred.Cells.Item(24, 4).Value = 0
$wsPred.Cells.Item(25, 4).Value = 0

# --- Continent_statistics: round Yield_* columns (C,D,E,F) to whole numbers ---
$wsCont.Cells.Item(3, 3).Value = 2085
$wsCont.Cells.Item(3, 4).Value = 2321
$wsCont.Cells.Item(3, 5).Value = 1594
$wsCont.Cells.Item(3, 6).Value = 1470
$wsCont.Cells.Item(4, 3).Value = 4656
$wsCont.Cells.Item(4, 4).Value = 5432
$wsCont.Cells.Item(4, 5).Value = 3301
$wsCont.Cells.Item(4, 6).Value = 2816
$wsCont.Cells.Item(5, 3).Value = 6414
$wsCont.Cells.Item(5, 4).Value = 6505
$wsCont.Cells.Item(5, 5).Value = 5271
$wsCont.Cells.Item(5, 6).Value = 3748
$wsCont.Cells.Item(6, 3).Value = 8597
$wsCont.Cells.Item(6, 4).Value = 6625
$wsCont.Cells.Item(6, 5).Value = 5708
$wsCont.Cells.Item(6, 6).Value = 3739
$wsCont.Cells.Item(7, 3).Value = 6369
$wsCont.Cells.Item(7, 4).Value = 5815
$wsCont.Cells.Item(7, 5).Value = 4974
$wsCont.Cells.Item(7, 6).Value = 3353
$wsCont.Cells.Item(8, 3).Value = 4526
$wsCont.Cells.Item(8, 4).Value = 4006
$wsCont.Cells.Item(8, 5).Value = 3208
$wsCont.Cells.Item(8, 6).Value = 2194
$wsCont.Cells.Item(9, 3).Value = 2650
$wsCont.Cells.Item(9, 4).Value = 2823
$wsCont.Cells.Item(9, 5).Value = 1985
$wsCont.Cells.Item(9, 6).Value = 1969
$wsCont.Cells.Item(10, 3).Value = 4425
$wsCont.Cells.Item(10, 4).Value = 4515
$wsCont.Cells.Item(10, 5).Value = 3021
$wsCont.Cells.Item(10, 6).Value = 2869
$wsCont.Cells.Item(11, 3).Value = 6722
$wsCont.Cells.Item(11, 4).Value = 5858
$wsCont.Cells.Item(11, 5).Value = 4770
$wsCont.Cells.Item(11, 6).Value = 4058
$wsCont.Cells.Item(12, 3).Value = 6391
$wsCont.Cells.Item(12, 4).Value = 5305
$wsCont.Cells.Item(12, 6).Value = 3650
$wsCont.Cells.Item(13, 3).Value = 8848
$wsCont.Cells.Item(13, 4).Value = 4844
$wsCont.Cells.Item(13, 5).Value = 4421
$wsCont.Cells.Item(13, 6).Value = 3757
$wsCont.Cells.Item(14, 3).Value = 5295
$wsCont.Cells.Item(14, 4).Value = 4380
$wsCont.Cells.Item(14, 5).Value = 3460
$wsCont.Cells.Item(14, 6).Value = 3053
$wsCont.Cells.Item(15, 3).Value = 999
$wsCont.Cells.Item(15, 4).Value = 1627
$wsCont.Cells.Item(15, 5).Value = 985
$wsCont.Cells.Item(15, 6).Value = 940
$wsCont.Cells.Item(16, 3).Value = 1424
$wsCont.Cells.Item(16, 4).Value = 1564
$wsCont.Cells.Item(16, 5).Value = 1189
$wsCont.Cells.Item(16, 6).Value = 1093
$wsCont.Cells.Item(17, 3).Value = 2523
$wsCont.Cells.Item(17, 4).Value = 2444
$wsCont.Cells.Item(17, 5).Value = 2205
$wsCont.Cells.Item(17, 6).Value = 1519
$wsCont.Cells.Item(18, 3).Value = 2916
$wsCont.Cells.Item(18, 4).Value = 2643
$wsCont.Cells.Item(18, 5).Value = 2468
$wsCont.Cells.Item(18, 6).Value = 1639
$wsCont.Cells.Item(19, 3).Value = 2685
$wsCont.Cells.Item(19, 4).Value = 2480
$wsCont.Cells.Item(19, 5).Value = 2289
$wsCont.Cells.Item(19, 6).Value = 1526
$wsCont.Cells.Item(20, 3).Value = 2546
$wsCont.Cells.Item(20, 4).Value = 2151
$wsCont.Cells.Item(20, 5).Value = 1408
$wsCont.Cells.Item(20, 6).Value = 1249
$wsCont.Cells.Item(21, 3).Value = 2873
$wsCont.Cells.Item(21, 4).Value = 3168
$wsCont.Cells.Item(21, 5).Value = 1807
$wsCont.Cells.Item(21, 6).Value = 1633
$wsCont.Cells.Item(22, 3).Value = 4040
$wsCont.Cells.Item(22, 4).Value = 3137
$wsCont.Cells.Item(22, 5).Value = 2336
$wsCont.Cells.Item(22, 6).Value = 1751
$wsCont.Cells.Item(23, 3).Value = 3021
$wsCont.Cells.Item(23, 4).Value = 2840
$wsCont.Cells.Item(23, 5).Value = 2110
$wsCont.Cells.Item(23, 6).Value = 1575
$wsCont.Cells.Item(24, 3).Value = 1762
$wsCont.Cells.Item(24, 4).Value = 1724
$wsCont.Cells.Item(24, 5).Value = 1418
$wsCont.Cells.Item(24, 6).Value = 1075
$wsCont.Cells.Item(25, 3).Value = 2819
$wsCont.Cells.Item(25, 4).Value = 3220
$wsCont.Cells.Item(25, 5).Value = 2219
$wsCont.Cells.Item(25, 6).Value = 1741
